$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the J:M block with new values (duplicated/refreshed data from F:I block)
$ws.Range("J1").Value = 3936
$ws.Range("K1").Value = 3937
$ws.Range("L1").Value = 3938
$ws.Range("M1").Value = 3939

$ws.Range("J2").Value = 39.6
$ws.Range("K2").Value = 42.4
$ws.Range("L2").Value = 33.4
$ws.Range("M2").Value = 41.6

$ws.Range("J3").Value = 39
$ws.Range("K3").Value = 41.9
$ws.Range("L3").Value = 33.1
$ws.Range("M3").Value = 41.3

$ws.Range("J4").Value = 41
$ws.Range("K4").Value = 42.1
$ws.Range("L4").Value = 34.5
$ws.Range("M4").Value = 43.3

$ws.Range("J5").Value = 39.1
$ws.Range("K5").Value = 41
$ws.Range("L5").Value = 34
$ws.Range("M5").Value = 42

# Drop the right-hand border on column I (it no longer divides two blocks since J:M was refreshed)
$ws.Range("I1:I5").Borders.Item(10).LineStyle = 0

# Move/restore the active selection to the refreshed J1:M1 header cells
$ws.Range("J1:M1").Select()
